$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-49): significance test name, test stat, p-value, and significance flag ---
$ws.Range("C2").Value = "One-way F test"
$ws.Range("D2").Value = [double]"103.9250343758279"
$ws.Range("E2").Value = [double]"5.450084891987035e-43"

$ws.Range("C3").Value = "One-way F test"
$ws.Range("D3").Value = [double]"492.2975528874661"
$ws.Range("E3").Value = [double]"4.982197960906943e-165"

$ws.Range("C4").Value = "One-way F test"
$ws.Range("D4").Value = [double]"2488.279675972683"
$ws.Range("E4").Value = [double]"0"

$ws.Range("C5").Value = "One-way F test"
$ws.Range("D5").Value = [double]"36.80602324941301"
$ws.Range("E5").Value = [double]"2.488443800449051e-16"

$ws.Range("C6").Value = "One-way F test"
$ws.Range("D6").Value = [double]"1865.500228308488"
$ws.Range("E6").Value = [double]"0"

$ws.Range("C7").Value = "One-way F test"
$ws.Range("D7").Value = [double]"340.6284512531364"
$ws.Range("E7").Value = [double]"1.199104472394488e-122"

$ws.Range("C8").Value = "One-way F test"
$ws.Range("D8").Value = [double]"93.14539548606869"
$ws.Range("E8").Value = [double]"7.471291394604005e-39"

$ws.Range("C9").Value = "One-way F test"
$ws.Range("D9").Value = [double]"17.23422415869616"
$ws.Range("E9").Value = [double]"3.982348758822768e-08"

$ws.Range("C10").Value = "One-way F test"
$ws.Range("D10").Value = [double]"4.10521155535292"
$ws.Range("E10").Value = [double]"0.01667249606816113"

$ws.Range("C11").Value = "One-way F test"
$ws.Range("D11").Value = [double]"24.28869629620745"
$ws.Range("E11").Value = [double]"4.159931931841053e-11"

$ws.Range("C12").Value = "One-way F test"
$ws.Range("D12").Value = [double]"259.0150279682046"
$ws.Range("E12").Value = [double]"2.504918173677693e-97"

$ws.Range("C13").Value = "One-way F test"
$ws.Range("D13").Value = [double]"292.7912155427557"
$ws.Range("E13").Value = [double]"4.783773674273814e-108"

$ws.Range("C14").Value = "T-Test"
$ws.Range("D14").Value = [double]"-3.20007098844747"
$ws.Range("E14").Value = [double]"0.001461579957827798"

$ws.Range("C15").Value = "T-Test"
$ws.Range("D15").Value = [double]"-4.62777566626418"
$ws.Range("E15").Value = [double]"4.72011467708349e-06"

$ws.Range("C16").Value = "T-Test"
$ws.Range("D16").Value = [double]"5.292175244158653"
$ws.Range("E16").Value = [double]"1.815062139272068e-07"

$ws.Range("C18").Value = "T-Test"
$ws.Range("D18").Value = [double]"-4.366461921739277"
$ws.Range("E18").Value = [double]"1.536775195372723e-05"

$ws.Range("C20").Value = "T-Test"
$ws.Range("D20").Value = [double]"2.897497209476683"
$ws.Range("E20").Value = [double]"0.003927216464984266"

$ws.Range("C22").Value = "T-Test"
$ws.Range("D22").Value = [double]"3.557182660562407"
$ws.Range("E22").Value = [double]"0.000410682646011469"

$ws.Range("C23").Value = "T-Test"
$ws.Range("D23").Value = [double]"-0.867142270016976"
$ws.Range("E23").Value = [double]"0.3862816838368289"

$ws.Range("C24").Value = "T-Test"
$ws.Range("D24").Value = [double]"-15.17561716272611"
$ws.Range("E24").Value = [double]"4.968122309413295e-43"

$ws.Range("C25").Value = "T-Test"
$ws.Range("D25").Value = [double]"15.41329739875774"
$ws.Range("E25").Value = [double]"4.144656743666704e-44"

$ws.Range("C26").Value = "T-Test"
$ws.Range("D26").Value = [double]"3.63490470761788"
$ws.Range("E26").Value = [double]"0.0003069922751638511"

$ws.Range("C27").Value = "T-Test"
$ws.Range("D27").Value = [double]"-3.879525504622102"
$ws.Range("E27").Value = [double]"0.0001187192107984613"

$ws.Range("C28").Value = "T-Test"
$ws.Range("D28").Value = [double]"3.703629135984921"
$ws.Range("E28").Value = [double]"0.000236309863528453"

$ws.Range("C29").Value = "T-Test"
$ws.Range("D29").Value = [double]"0.5409203513035115"
$ws.Range("E29").Value = [double]"0.5888043460659246"

$ws.Range("C30").Value = "T-Test"
$ws.Range("D30").Value = [double]"-3.226509327587789"
$ws.Range("E30").Value = [double]"0.001335558424313903"

$ws.Range("C31").Value = "T-Test"
$ws.Range("D31").Value = [double]"3.789572000834149"
$ws.Range("E31").Value = [double]"0.0001693798411145114"

$ws.Range("C32").Value = "T-Test"
$ws.Range("D32").Value = [double]"-0.2054382678130561"
$ws.Range("E32").Value = [double]"0.8373136719515815"

$ws.Range("C33").Value = "T-Test"
$ws.Range("D33").Value = [double]"0.4464872046650555"
$ws.Range("E33").Value = [double]"0.6554394366694096"

$ws.Range("C34").Value = "T-Test"
$ws.Range("D34").Value = [double]"0.6760445824019018"
$ws.Range("E34").Value = [double]"0.4993261673274386"

$ws.Range("C35").Value = "T-Test"
$ws.Range("D35").Value = [double]"1.2416706355686"
$ws.Range("E35").Value = [double]"0.2149425882217344"

$ws.Range("C36").Value = "T-Test"
$ws.Range("D36").Value = [double]"1.967883838387863"
$ws.Range("E36").Value = [double]"0.04963582312876014"
$ws.Range("F36").Value = $true

$ws.Range("C37").Value = "T-Test"
$ws.Range("D37").Value = [double]"-1.210112196764809"
$ws.Range("E37").Value = [double]"0.22680991195355"

$ws.Range("C38").Value = "T-Test"
$ws.Range("D38").Value = [double]"-26.70903432385005"
$ws.Range("E38").Value = [double]"3.478876148824275e-98"

$ws.Range("C39").Value = "T-Test"
$ws.Range("D39").Value = [double]"22.40540323826895"
$ws.Range("E39").Value = [double]"2.050111511544678e-77"

$ws.Range("C40").Value = "T-Test"
$ws.Range("D40").Value = [double]"18.07577900150219"
$ws.Range("E40").Value = [double]"1.587913225899474e-56"

$ws.Range("C41").Value = "T-Test"
$ws.Range("D41").Value = [double]"-20.3533957131679"
$ws.Range("E41").Value = [double]"1.836427519341772e-67"

$ws.Range("C42").Value = "T-Test"
$ws.Range("D42").Value = [double]"23.35243727784407"
$ws.Range("E42").Value = [double]"5.19572498798019e-82"

$ws.Range("C43").Value = "T-Test"
$ws.Range("D43").Value = [double]"-24.28682441786492"
$ws.Range("E43").Value = [double]"1.547847087079104e-86"

$ws.Range("C44").Value = "T-Test"
$ws.Range("D44").Value = [double]"-2.016040750719476"
$ws.Range("E44").Value = [double]"0.04433218659223444"

$ws.Range("C45").Value = "T-Test"
$ws.Range("D45").Value = [double]"2.659495576454042"
$ws.Range("E45").Value = [double]"0.008077905587840311"

$ws.Range("C46").Value = "T-Test"
$ws.Range("D46").Value = [double]"-0.8638296036201208"
$ws.Range("E46").Value = [double]"0.3880973580777607"

$ws.Range("C47").Value = "T-Test"
$ws.Range("D47").Value = [double]"-1.921300344464408"
$ws.Range("E47").Value = [double]"0.05526432420963388"
$ws.Range("F47").Value = $false

$ws.Range("C48").Value = "T-Test"
$ws.Range("D48").Value = [double]"-13.99724859110737"
$ws.Range("E48").Value = [double]"8.894914808069599e-38"

$ws.Range("C49").Value = "T-Test"
$ws.Range("D49").Value = [double]"13.79218912698009"
$ws.Range("E49").Value = [double]"7.004891698637221e-37"

# --- Append new rows (50-61) for the IAS comparison group ---
$ws.Range("A50").Value = "IAS"
$ws.Range("B50").Value = "ratio_paras"
$ws.Range("C50").Value = "One-way F test"
$ws.Range("D50").Value = [double]"337.1931345863907"
$ws.Range("E50").Value = [double]"5.072061204772169e-242"
$ws.Range("F50").Value = $true

$ws.Range("A51").Value = "IAS"
$ws.Range("B51").Value = "ratio_list_items"
$ws.Range("C51").Value = "One-way F test"
$ws.Range("D51").Value = [double]"503.0133287261479"
$ws.Range("E51").Value = [double]"4.474863225088838e-317"
$ws.Range("F51").Value = $true

$ws.Range("A52").Value = "IAS"
$ws.Range("B52").Value = "ratio_headings"
$ws.Range("C52").Value = "One-way F test"
$ws.Range("D52").Value = [double]"1434.350177387523"
$ws.Range("E52").Value = [double]"0"
$ws.Range("F52").Value = $true

$ws.Range("A53").Value = "IAS"
$ws.Range("B53").Value = "avg_para_len"
$ws.Range("C53").Value = "One-way F test"
$ws.Range("D53").Value = [double]"121.9227978430191"
$ws.Range("E53").Value = [double]"2.327858259410091e-108"
$ws.Range("F53").Value = $true

$ws.Range("A54").Value = "IAS"
$ws.Range("B54").Value = "num_sentences"
$ws.Range("C54").Value = "One-way F test"
$ws.Range("D54").Value = [double]"854.5093066193562"
$ws.Range("E54").Value = [double]"0"
$ws.Range("F54").Value = $true

$ws.Range("A55").Value = "IAS"
$ws.Range("B55").Value = "avg_len"
$ws.Range("C55").Value = "One-way F test"
$ws.Range("D55").Value = [double]"600.8964257179297"
$ws.Range("E55").Value = [double]"0"
$ws.Range("F55").Value = $true

$ws.Range("A56").Value = "IAS"
$ws.Range("B56").Value = "flesch"
$ws.Range("C56").Value = "One-way F test"
$ws.Range("D56").Value = [double]"40.10816492550137"
$ws.Range("E56").Value = [double]"8.742930203466194e-39"
$ws.Range("F56").Value = $true

$ws.Range("A57").Value = "IAS"
$ws.Range("B57").Value = "cli"
$ws.Range("C57").Value = "One-way F test"
$ws.Range("D57").Value = [double]"11.30245280737467"
$ws.Range("E57").Value = [double]"1.017809412477035e-10"
$ws.Range("F57").Value = $true

$ws.Range("A58").Value = "IAS"
$ws.Range("B58").Value = "avg_concrete"
$ws.Range("C58").Value = "One-way F test"
$ws.Range("D58").Value = [double]"2.65397987484021"
$ws.Range("E58").Value = [double]"0.02136283399874483"
$ws.Range("F58").Value = $true

$ws.Range("A59").Value = "IAS"
$ws.Range("B59").Value = "concrete_ratio"
$ws.Range("C59").Value = "One-way F test"
$ws.Range("D59").Value = [double]"11.670469474077"
$ws.Range("E59").Value = [double]"4.38810462932117e-11"
$ws.Range("F59").Value = $true

$ws.Range("A60").Value = "IAS"
$ws.Range("B60").Value = "abstract_ratio"
$ws.Range("C60").Value = "One-way F test"
$ws.Range("D60").Value = [double]"221.1192272393743"
$ws.Range("E60").Value = [double]"8.646651848982935e-177"
$ws.Range("F60").Value = $true

$ws.Range("A61").Value = "IAS"
$ws.Range("B61").Value = "undefined_ratio"
$ws.Range("C61").Value = "One-way F test"
$ws.Range("D61").Value = [double]"228.4717117168348"
$ws.Range("E61").Value = [double]"2.473556392891394e-181"
$ws.Range("F61").Value = $true
